$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 7
$ws.Range("F4").Value = 1868
$ws.Range("F5").Value = 3258
$ws.Range("F6").Value = 204
$ws.Range("F7").Value = 4766
$ws.Range("F8").Value = 451
$ws.Range("F9").Value = 278
$ws.Range("F10").Value = 161
$ws.Range("F11").Value = 617
$ws.Range("F13").Value = 19
$ws.Range("F14").Value = 9
$ws.Range("F15").Value = 645
$ws.Range("F17").Value = 18
$ws.Range("F20").Value = 333
$ws.Range("F21").Value = 4708
$ws.Range("F22").Value = 14
$ws.Range("F23").Value = 33
$ws.Range("F25").Value = 5853
$ws.Range("F26").Value = 13
$ws.Range("F27").Value = 1187
$ws.Range("F28").Value = 237
$ws.Range("F29").Value = 656
$ws.Range("F30").Value = 4412
$ws.Range("F31").Value = 4
$ws.Range("F32").Value = 81
$ws.Range("F33").Value = 120
$ws.Range("F34").Value = 821
$ws.Range("F35").Value = 62
$ws.Range("F36").Value = 755
$ws.Range("F37").Value = 781

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 34
$ws.Range("F6").Value = 45

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 7
$ws.Range("F7").Value = 1868
$ws.Range("F9").Value = 3258
$ws.Range("F10").Value = 204
$ws.Range("F11").Value = 4766
$ws.Range("F12").Value = 451
$ws.Range("F13").Value = 278
$ws.Range("F14").Value = 161
$ws.Range("F15").Value = 617
$ws.Range("F17").Value = 19
$ws.Range("F18").Value = 9
$ws.Range("F19").Value = 645
$ws.Range("F21").Value = 18
$ws.Range("F22").Value = 34
$ws.Range("F25").Value = 333
$ws.Range("F26").Value = 4708
$ws.Range("F27").Value = 14
$ws.Range("F28").Value = 33
$ws.Range("F30").Value = 5853
$ws.Range("F31").Value = 13
$ws.Range("F32").Value = 1187
$ws.Range("F33").Value = 237
$ws.Range("F34").Value = 656
$ws.Range("F35").Value = 4412
$ws.Range("F36").Value = 4
$ws.Range("F38").Value = 81
$ws.Range("F39").Value = 120
$ws.Range("F40").Value = 821
$ws.Range("F41").Value = 62
$ws.Range("F42").Value = 755
$ws.Range("F43").Value = 781
$ws.Range("F45").Value = 45
